# Auto-generated update of cached market-price snapshot cells
# (per-cell values refreshed by the scheduled Hyperion Profits runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1419.4166
$ws.Range("J2").Value = 1083
$ws.Range("L2").Value = 1083
$ws.Range("N2").Value = -1309
$ws.Range("H132").Value = 23258792
$ws.Range("I132").Value = 25644124
$ws.Range("K132").Value = 76932372
$ws.Range("M132").Value = -76929842
$ws.Range("H135").Value = 1847.4736
$ws.Range("I135").Value = 901.8889
$ws.Range("K135").Value = 8117.0001
$ws.Range("M135").Value = -5582.0001
$ws.Range("H138").Value = 2730.4736
$ws.Range("I138").Value = 1240.6666
$ws.Range("J138").Value = 3418.077
$ws.Range("K138").Value = 3721.9998
$ws.Range("L138").Value = 10254.231
$ws.Range("M138").Value = 1418.0002
$ws.Range("N138").Value = -20534.231
$ws.Range("H141").Value = 3155.4546
$ws.Range("I141").Value = 3145
$ws.Range("J141").Value = 3183.3333
$ws.Range("K141").Value = 9435
$ws.Range("L141").Value = 9549.999899999999
$ws.Range("M141").Value = -4255
$ws.Range("N141").Value = -19909.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6542003.5
$ws.Range("I45").Value = 11067322
$ws.Range("K45").Value = 11067322
$ws.Range("M45").Value = -11066945
$ws.Range("H61").Value = 3591.7856
$ws.Range("I61").Value = 2999
$ws.Range("K61").Value = 2999
$ws.Range("M61").Value = -2787
$ws.Range("H132").Value = 3067.0908
$ws.Range("I132").Value = 2182.4
$ws.Range("K132").Value = 6547.200000000001
$ws.Range("M132").Value = -4017.200000000001
$ws.Range("H136").Value = 3591.7856
$ws.Range("I136").Value = 2999
$ws.Range("K136").Value = 8997
$ws.Range("M136").Value = -6447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1448.1538
$ws.Range("I20").Value = 1414.8182
$ws.Range("K20").Value = 1414.8182
$ws.Range("M20").Value = -1167.8182
$ws.Range("H86").Value = 5564107
$ws.Range("I86").Value = 5891261
$ws.Range("J86").Value = 2495
$ws.Range("K86").Value = 5891261
$ws.Range("L86").Value = 2495
$ws.Range("M86").Value = -5890138
$ws.Range("N86").Value = -4741
$ws.Range("H89").Value = 5564107
$ws.Range("I89").Value = 5891261
$ws.Range("J89").Value = 2495
$ws.Range("K89").Value = 29456305
$ws.Range("L89").Value = 12475
$ws.Range("M89").Value = -29450689
$ws.Range("N89").Value = -23707
$ws.Range("H97").Value = 4955.5
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 296.14285
$ws.Range("I7").Value = 104.454544
$ws.Range("K7").Value = 104.454544
$ws.Range("M7").Value = 8.545456000000001
$ws.Range("H8").Value = 450.5
$ws.Range("J8").Value = 450.5
$ws.Range("L8").Value = 450.5
$ws.Range("N8").Value = -730.5
$ws.Range("H31").Value = 30493.156
$ws.Range("I31").Value = 1439.75
$ws.Range("K31").Value = 1439.75
$ws.Range("M31").Value = -1144.75
$ws.Range("H34").Value = 30493.156
$ws.Range("I34").Value = 1439.75
$ws.Range("K34").Value = 1439.75
$ws.Range("M34").Value = -1237.75
$ws.Range("H58").Value = 2052.7144
$ws.Range("I58").Value = 1710.6471
$ws.Range("K58").Value = 1710.6471
$ws.Range("M58").Value = -1507.6471
$ws.Range("H134").Value = 49069.05
$ws.Range("I134").Value = 115435.75
$ws.Range("J134").Value = 4824.5835
$ws.Range("K134").Value = 346307.25
$ws.Range("L134").Value = 14473.7505
$ws.Range("M134").Value = -343772.25
$ws.Range("N134").Value = -19543.7505
$ws.Range("H136").Value = 2052.7144
$ws.Range("I136").Value = 1710.6471
$ws.Range("K136").Value = 5131.9413
$ws.Range("M136").Value = -2581.9413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3000.5715
$ws.Range("I69").Value = 2001.3334
$ws.Range("J69").Value = 3750
$ws.Range("K69").Value = 6004.0002
$ws.Range("L69").Value = 11250
$ws.Range("M69").Value = -5193.0002
$ws.Range("N69").Value = -12872
$ws.Range("H72").Value = 3000.5715
$ws.Range("I72").Value = 2001.3334
$ws.Range("J72").Value = 3750
$ws.Range("K72").Value = 18012.0006
$ws.Range("L72").Value = 33750
$ws.Range("M72").Value = -13956.0006
$ws.Range("N72").Value = -41862
$ws.Range("H121").Value = 635
$ws.Range("I121").Value = 235.625
$ws.Range("J121").Value = 2232.5
$ws.Range("K121").Value = 706.875
$ws.Range("L121").Value = 6697.5
$ws.Range("M121").Value = 603.125
$ws.Range("N121").Value = -9317.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 21599.8
$ws.Range("I58").Value = 9333
$ws.Range("K58").Value = 9333
$ws.Range("M58").Value = -9056
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 406853
$ws.Range("J122").Value = 5497.5
$ws.Range("L122").Value = 16492.5
$ws.Range("N122").Value = -21392.5
$ws.Range("H126").Value = 5369433.5
$ws.Range("I126").Value = 3249517.8
$ws.Range("J126").Value = 8337315.5
$ws.Range("K126").Value = 9748553.399999999
$ws.Range("L126").Value = 25011946.5
$ws.Range("M126").Value = -9746083.399999999
$ws.Range("N126").Value = -25016886.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3475566.2
$ws.Range("J61").Value = 2195.7693
$ws.Range("L61").Value = 2195.7693
$ws.Range("N61").Value = -2599.7693
$ws.Range("H93").Value = 19609538
$ws.Range("I93").Value = 27779652
$ws.Range("K93").Value = 27779652
$ws.Range("M93").Value = -27778404
$ws.Range("H100").Value = 3298
$ws.Range("I100").Value = 3076.6924
$ws.Range("J100").Value = 3777.5
$ws.Range("K100").Value = 3076.6924
$ws.Range("L100").Value = 3777.5
$ws.Range("M100").Value = -2535.6924
$ws.Range("N100").Value = -4859.5
$ws.Range("H113").Value = 3475566.2
$ws.Range("J113").Value = 2195.7693
$ws.Range("L113").Value = 2195.7693
$ws.Range("N113").Value = -6535.7693
$ws.Range("I122").Value = 4580.6
$ws.Range("J122").Value = 6680.727
$ws.Range("K122").Value = 13741.8
$ws.Range("L122").Value = 20042.181
$ws.Range("M122").Value = -11291.8
$ws.Range("N122").Value = -24942.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11112064
$ws.Range("I81").Value = 15152488
$ws.Range("K81").Value = 30304976
$ws.Range("M81").Value = -30303915
$ws.Range("H84").Value = 11112064
$ws.Range("I84").Value = 15152488
$ws.Range("K84").Value = 151524880
$ws.Range("M84").Value = -151519576
$ws.Range("H136").Value = 1536.1111
$ws.Range("I136").Value = 1119
$ws.Range("K136").Value = 3357
$ws.Range("M136").Value = -807
